{"js": "// Capitalize \"rasi\" -> \"Rasi\" in the phrase \"rasi bintang Sagitarius\"\n// (e.g. \"Waktu Kampanye 2022 untuk rasi bintang Sagitarius\" and\n// \"... identifikasi  rasi bintang Sagitarius ...\") throughout the document,\n// while leaving unrelated occurrences of the standalone word \"rasi\"\n// (e.g. \"rasi bintang yang dituju\") untouched.\n\nconst body = context.document.body;\nconst results = body.search(\"rasi bintang Sagitarius\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nresults.items.forEach((range) => {\n  range.insertText(\"Rasi bintang Sagitarius\", Word.InsertLocation.replace);\n});\n\nawait context.sync();\n", "ps1": "# Capitalize \"rasi\" -> \"Rasi\" in the phrase \"rasi bintang Sagitarius\"\n# (e.g. \"Waktu Kampanye 2022 untuk rasi bintang Sagitarius\" and\n# \"... identifikasi  rasi bintang Sagitarius ...\") throughout the document,\n# while leaving unrelated occurrences of the standalone word \"rasi\"\n# (e.g. \"rasi bintang yang dituju\") untouched.\n\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"rasi bintang Sagitarius\", $true, $false, $false, $false, $false, $true, 1, $false, \"Rasi bintang Sagitarius\", 2)\n"}
